$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = '36.433.58'
$ws.Cells.Item(2,5).Value = '  -0.17%  '

$ws.Cells.Item(3,4).Value = '1.938.30'
$ws.Cells.Item(3,5).Value = '  -1.29%  '

$ws.Cells.Item(4,5).Value = '  -0.01%  '

$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = '242.35'
$ws.Cells.Item(5,5).Value = '  -0.77%  '

$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = '0.608'
$ws.Cells.Item(6,5).Value = '  -1.76%  '

$ws.Cells.Item(7,5).Value = '  -0.09%  '

$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = '56.38'
$ws.Cells.Item(8,5).Value = '  -4.29%  '

$ws.Cells.Item(9,5).Value = '  -4.24%  '

$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = '0.0803'
$ws.Cells.Item(10,5).Value = '  -5.45%  '

$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = '0.102'
$ws.Cells.Item(11,5).Value = '  -1.17%  '

$ws.Cells.Item(12,4).Value = '2.217.81'
$ws.Cells.Item(12,5).Value = '  -1.55%  '

$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = '21.23'
$ws.Cells.Item(13,5).Value = '  -3.91%  '

$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = '0.801'
$ws.Cells.Item(14,5).Value = '  -4.77%  '

$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = '13.26'
$ws.Cells.Item(15,5).Value = '  -2.36%  '

$ws.Cells.Item(16,5).Value = '  -4.39%  '

$ws.Cells.Item(17,4).Value = '1.936.49'
$ws.Cells.Item(17,5).Value = '  -1.85%  '

$ws.Cells.Item(18,4).Value = '36.408.12'
$ws.Cells.Item(18,5).Value = '  -0.04%  '

$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = '68.81'
$ws.Cells.Item(19,5).Value = '  -2.23%  '

$ws.Cells.Item(20,4).Value = '0.0₃0851'
$ws.Cells.Item(20,5).Value = '  -3.89%  '

$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = '226.53'
$ws.Cells.Item(21,5).Value = '  -2.15%  '

$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = '4.92'
$ws.Cells.Item(22,5).Value = '  -3.36%  '

$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = '1.00'
$ws.Cells.Item(23,5).Value = '  +0.14%  '

$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = '2.40'
$ws.Cells.Item(24,5).Value = '  -5.34%  '

$ws.Cells.Item(25,5).Value = '  +0.05%  '

$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = '9.06'
$ws.Cells.Item(26,5).Value = '  -5.21%  '

$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = '159.60'
$ws.Cells.Item(27,5).Value = '  -3.28%  '

$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value = '0.132'
$ws.Cells.Item(28,5).Value = '  +8.19%  '

$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = '19.02'
$ws.Cells.Item(29,5).Value = '  -3.70%  '

$ws.Cells.Item(30,5).Value = '  -1.24%  '

$ws.Cells.Item(31,5).Value = '  -7.31%  '

$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(32,4).Value = '4.54'
$ws.Cells.Item(32,5).Value = '  -4.85%  '

$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = '0.0611'
$ws.Cells.Item(33,5).Value = '  -4.43%  '

$ws.Cells.Item(34,2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(34,3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = '4.11'
$ws.Cells.Item(34,5).Value = '  -5.75%  '

$ws.Cells.Item(35,2).Value = 'BinanceUSD'
$ws.Cells.Item(35,3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value = '1.00'
$ws.Cells.Item(35,5).Value = '  -0.07%  '

$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value = '6.09'
$ws.Cells.Item(36,5).Value = '  -1.32%  '

$ws.Cells.Item(37,5).Value = '  -1.26%  '

$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = '2.18'
$ws.Cells.Item(38,5).Value = '  +0.30%  '

$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value = '3.17'
$ws.Cells.Item(39,5).Value = '  +9.03%  '

$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = '0.0988'
$ws.Cells.Item(40,5).Value = '  +0.32%  '

$ws.Cells.Item(41,5).Value = '  +1.16%  '

$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = '0.0208'
$ws.Cells.Item(42,5).Value = '  -1.49%  '

$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = '1.14'
$ws.Cells.Item(43,5).Value = '  -4.49%  '

$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = '15.61'
$ws.Cells.Item(44,5).Value = '  -1.62%  '

$ws.Cells.Item(45,2).Value = 'ARBITRUM'
$ws.Cells.Item(45,3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = '1.02'
$ws.Cells.Item(45,5).Value = '  -3.13%  '

$ws.Cells.Item(46,2).Value = 'Maker'
$ws.Cells.Item(46,3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(46,4).Value = '1.328.32'
$ws.Cells.Item(46,5).Value = '  -1.56%  '

$ws.Cells.Item(47,5).Value = '  -4.19%  '

$ws.Cells.Item(48,5).Value = '  -4.07%  '

$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = '2.81'
$ws.Cells.Item(49,5).Value = '  -0.60%  '

$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = '3.50'
$ws.Cells.Item(50,5).Value = '  +14.31%  '

$ws.Cells.Item(51,4).Value = '2.111.48'
$ws.Cells.Item(51,5).Value = '  -1.53%  '
